$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - "Learn networking fundamentals: IP, DNS, TCP/UDP" (Day 17)
# Status: In Progress -> Done
$ws.Range("C8").Value = "Done"
# In Progress? flag: checked -> unchecked
$ws.Range("D8").Value = "☐"
# Done? flag: unchecked -> checked
$ws.Range("E8").Value = "☑"

# Row 9 - "Security foundations: threat vs vulnerability, risk concepts" (Day 19)
# Status: Not Started -> In Progress
$ws.Range("C9").Value = "In Progress"
# In Progress? flag: unchecked -> checked
$ws.Range("D9").Value = "☑"

# Update selection to reflect last edited cell
$ws.Range("D9").Select()
